$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.224.00"
$ws.Range("E2").Value = "  -1.69%  "
$ws.Range("D3").Value = "1.877.40"
$ws.Range("E3").Value = "  -0.75%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "236.36"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.95%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9998"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.05%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4840"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -1.11%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2875"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -3.28%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06591"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -2.21%  "
$ws.Range("D10").Value = "1.887.14"
$ws.Range("E10").Value = "  +0.95%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "16.71"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -1.00%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07282"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -0.31%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.182"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +1.50%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "87.10"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -2.87%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6544"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -1.86%  "
$ws.Range("D16").Value = "30.195.60"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.35"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.68%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9995"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.12%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007698"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -2.98%  "
$ws.Range("D20").Value = "2.111.04"
$ws.Range("E20").Value = "  -0.21%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.313"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +7.66%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.002"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.00%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "195.75"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -5.67%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.123"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -1.26%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.310"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -3.20%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "158.95"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.60%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.07"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -4.25%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.916"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +2.03%  "
$ws.Range("E29").Value = "  +1.75%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.271"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -1.36%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09130"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.03%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.060"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.51%  "
$ws.Range("E33").Value = "  -0.84%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7195"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -3.91%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.095"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.93%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.714"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +1.11%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01796"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -1.84%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.637"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -2.29%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.9170"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.81%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.040"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -2.70%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "106.14"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.44%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4285"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -4.33%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.794"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.02%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9981"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.15%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "66.30"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.71%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.398"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -4.62%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1319"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -3.34%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.168"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +3.69%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3825"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -5.87%  "

# Row 49/50: Cronos <-> Elrond swap with updated values
$ws.Range("B49").Value = "Elrond"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "34.00"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -2.36%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05751"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -2.53%  "
